$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("哈希")

# Copy the style (s="4", wrap text) from the row above into the new cells C14:G14
$ws.Range("C13:G13").Copy()
$ws.Range("C14:G14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in row 14 -- LeetCode #219, "Contains Duplicate II" (hash map solution)
$ws.Range("B14").Value = 219
$ws.Range("C14").Value = "给定一个整数数组和一个整数 k，判断数组中是否存在两个不同的索引 i 和 j，使 nums [i] = nums [j]，并且 i 和 j 的差的 绝对值至多为 k。 `n 示例 1:`n 输入: nums = [1,2,3,1], k = 3`n输出: true "
$ws.Range("E14").Value = "哈希表`n索引"
$ws.Range("D14").Value = "1 迭代数组，存入map，key是元素值，value是其索引值`n2 判断元素是否在map中，如果不在就加入map`n3 如果存在`n    判断其索引值差值<=k，return true`n    判断其索引值差值>k，将现在的index索引值保存，继续迭代`n4 迭代结束条件是遍历完数组元素`n另一方案是使用集合，基于长度是k的窗口：`n         1 判断集合中是否存在此元素，不存在就添加此元素`n         2 如果集合元素个数超过k，将最前面元素删除，始终保持集合元素个数是k`n         3 如果集合中包含此元素，就是满足需要"
$ws.Range("F14").Value = "O(N), N是元素个数"
$ws.Range("G14").Value = "O(N)"

# Row needs to grow to fit the long wrapped text
$ws.Rows.Item(14).RowHeight = 308

# Selection left on the newly-filled method cell
$ws.Activate()
$ws.Range("D14").Select()
